$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 1 and 2 (кола/90 now first, вода/80 now second) ---
$ws.Range("A1").Value = "кола"
$ws.Range("B1").Value = 90
$ws.Range("A2").Value = "вода"
$ws.Range("B2").Value = 80

# --- Row 3: update price for чипсы (74 -> 40) ---
$ws.Range("B3").Value = 40

# --- Rows 4 and 5 (сок / сэндвич) are unchanged ---

# --- Add new rows 6-8, copying formatting from row 5 then filling values ---
$ws.Range("A5:B5").Copy()
$ws.Range("A6:B6").PasteSpecial(-4122)
$ws.Range("A7:B7").PasteSpecial(-4122)
$ws.Range("A8:B8").PasteSpecial(-4122)

$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 30

$ws.Range("A6").Value = "капучино"
$ws.Range("B6").Value = 120

$ws.Range("A7").Value = "латте"
$ws.Range("B7").Value = 120

$ws.Range("A8").Value = "раф"
$ws.Range("B8").Value = 240

# --- Update selection to reflect the newly active cell ---
$ws.Range("B8").Select()
